$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-theme the deck: the custom "Integral" theme colours are replaced with
#    the built-in Office Theme palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# ---------------------------------------------------------------------------
function Set-ThemeColor {
    param($scheme, [int]$index, [string]$rrggbb)
    $r = [Convert]::ToInt32($rrggbb.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($rrggbb.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($rrggbb.Substring(4, 2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $scheme.Colors($index).RGB = $bgr
}

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

Set-ThemeColor $themeColors 1  "000000"  # dk1
Set-ThemeColor $themeColors 2  "FFFFFF"  # lt1
Set-ThemeColor $themeColors 3  "44546A"  # dk2
Set-ThemeColor $themeColors 4  "E7E6E6"  # lt2
Set-ThemeColor $themeColors 5  "5B9BD5"  # accent1
Set-ThemeColor $themeColors 6  "ED7D31"  # accent2
Set-ThemeColor $themeColors 7  "A5A5A5"  # accent3
Set-ThemeColor $themeColors 8  "FFC000"  # accent4
Set-ThemeColor $themeColors 9  "4472C4"  # accent5
Set-ThemeColor $themeColors 10 "70AD47"  # accent6
Set-ThemeColor $themeColors 11 "0563C1"  # hyperlink
Set-ThemeColor $themeColors 12 "954F72"  # followed hyperlink

# ---------------------------------------------------------------------------
# 2) Update the table style used by the "PLENARY" table on slide 16 to the
#    style applied together with the new theme.
# ---------------------------------------------------------------------------
foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq "{51713AF5-90D0-47E0-9274-AE64739708DE}") {
                $table.ApplyStyle("{6FE67144-4486-4918-9FE5-030BB7D35592}")
            }
        }
    }
}
